$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 0.0495
$ws.Range("E2").Value = -0.21
$ws.Range("G2").Value = 0.1440173253925284
$ws.Range("H2").Value = 0.1440173253925284
$ws.Range("I2").Value = 0.04206821873308067
$ws.Range("J2").Value = 0.04130698429886302
$ws.Range("K2").Value = 5.16
$ws.Range("L2").Value = 0.02793719545208447
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("U2").Value = 4.76
$ws.Range("V2").Value = 0.02218080149114632
$ws.Range("W2").Value = 0.04414029084687767
$ws.Range("X2").Value = 0.04917661515342803
$ws.Range("Y2").Value = -0.005036324306550356
$ws.Range("Z2").Value = 1.413159908186687
$ws.Range("AA2").Value = 0.05837337413925019
$ws.Range("AB2").Value = 0.04610843143785517
$ws.Range("AC2").Value = 0.01226494270139502
$ws.Range("AD2").Value = 26.6
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 26.6
$ws.Range("AG2").Value = 21.84
$ws.Range("AH2").Value = 0.1102819237147595
$ws.Range("AI2").Value = 0.1953010279001469
$ws.Range("AJ2").Value = 0.09237015733378448
$ws.Range("AK2").Value = 0.1661594643944005
$ws.Range("AL2").Value = 5.25
$ws.Range("AM2").Value = 5.25
$ws.Range("AN2").Value = 2.703252032520326
$ws.Range("AO2").Value = 1.48
$ws.Range("AP2").Value = 2.219512195121951
$ws.Range("AQ2").Value = 1.48

# Row 3 updates
$ws.Range("D3").Value = 0.0495
$ws.Range("E3").Value = -0.21
$ws.Range("G3").Value = 0.1440173253925284
$ws.Range("H3").Value = 0.1440173253925284
$ws.Range("I3").Value = 0.04206821873308067
$ws.Range("J3").Value = 0.04130698429886302
$ws.Range("K3").Value = 5.16
$ws.Range("L3").Value = 0.02793719545208447
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 4.76
$ws.Range("V3").Value = 0.02218080149114632
$ws.Range("W3").Value = 0.04414029084687767
$ws.Range("X3").Value = 0.04917661515342803
$ws.Range("Y3").Value = -0.005036324306550356
$ws.Range("Z3").Value = 1.413159908186687
$ws.Range("AA3").Value = 0.05837337413925019
$ws.Range("AB3").Value = 0.04610843143785517
$ws.Range("AC3").Value = 0.01226494270139502
$ws.Range("AD3").Value = 26.6
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 26.6
$ws.Range("AG3").Value = 21.84
$ws.Range("AH3").Value = 0.1102819237147595
$ws.Range("AI3").Value = 0.1953010279001469
$ws.Range("AJ3").Value = 0.09237015733378448
$ws.Range("AK3").Value = 0.1661594643944005
$ws.Range("AL3").Value = 5.25
$ws.Range("AM3").Value = 5.25
$ws.Range("AN3").Value = 2.703252032520326
$ws.Range("AO3").Value = 1.48
$ws.Range("AP3").Value = 2.219512195121951
$ws.Range("AQ3").Value = 1.48

# Remove buybacks_cash_returned column (T) for both rows
$ws.Range("T2").ClearContents()
$ws.Range("T3").ClearContents()

